$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the two oldest data rows (2008年, 2009年). This shifts the
# 2010年/2011年 rows up to become rows 2 and 3, matching the target
# workbook (which now only spans rows 1-3 / A1:AT3).
$ws.Rows("2:3").Delete()
